# Weekly CompStat report roll-forward: "New crime data collected"
#
# 1) Bump the report's "Volume .. Number .." and "Report Covering the Week .."
#    rich-text captions by one week (report index 46 -> 47; week of
#    11/11/2024-11/17/2024 -> 11/18/2024-11/24/2024). These shared strings are
#    built from several same-formatted runs, so only the specific digit/date
#    substrings are replaced in place via Characters(start,len).Text, leaving
#    the rest of the caption untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Characters(21, 2).Text = "47"
$ws.Range("C9").Characters(27, 10).Text = "11/18/2024"
$ws.Range("C9").Characters(48, 10).Text = "11/24/2024"

# 2) Refresh the weekly/28-day/YTD/2-year crime-count table (rows 14-33) with
#    the newly collected figures. A handful of cells were previously blank,
#    shown via the shared placeholder text "0" (style 13, General format);
#    now that they have real counts, give them the same "#,##0" number
#    format the rest of the count columns use (style 14) so they read as
#    numbers instead of literal text.
foreach ($addr in @("C15", "F15", "C27", "F27", "C33", "F33")) {
    $ws.Range($addr).NumberFormat = "#,##0"
}

$ws.Range("D14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 21
$ws.Range("K14").Value = -47.619047619047
$ws.Range("M14").Value = -63.333333333333
$ws.Range("N14").Value = -90
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 70
$ws.Range("K15").Value = 37.254901960784
$ws.Range("L15").Value = 42.857142857142
$ws.Range("M15").Value = 18.64406779661
$ws.Range("N15").Value = -40.17094017094
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -28.571428571428
$ws.Range("F16").Value = 57
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 586
$ws.Range("J16").Value = 628
$ws.Range("K16").Value = -6.687898089171
$ws.Range("L16").Value = -16.045845272206
$ws.Range("M16").Value = -14.327485380117
$ws.Range("N16").Value = -79.219858156028
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 82
$ws.Range("G17").Value = 101
$ws.Range("H17").Value = -18.811881188118
$ws.Range("I17").Value = 1002
$ws.Range("J17").Value = 990
$ws.Range("K17").Value = 1.212121212121
$ws.Range("L17").Value = 0.50150451354
$ws.Range("M17").Value = 44.797687861271
$ws.Range("N17").Value = -26.431718061674
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 27.777777777777
$ws.Range("I18").Value = 262
$ws.Range("J18").Value = 320
$ws.Range("K18").Value = -18.125
$ws.Range("L18").Value = -24.71264367816
$ws.Range("M18").Value = -40.589569160997
$ws.Range("N18").Value = -84.414039262343
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -31.578947368421
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -15.189873417721
$ws.Range("I19").Value = 820
$ws.Range("J19").Value = 985
$ws.Range("K19").Value = -16.751269035533
$ws.Range("L19").Value = -31.495405179615
$ws.Range("M19").Value = 28.125
$ws.Range("N19").Value = -4.982618771726
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 42
$ws.Range("G20").Value = 53
$ws.Range("H20").Value = -20.754716981132
$ws.Range("I20").Value = 548
$ws.Range("J20").Value = 500
$ws.Range("K20").Value = 9.6
$ws.Range("L20").Value = 19.130434782608
$ws.Range("M20").Value = 93.639575971731
$ws.Range("N20").Value = -76.955424726661
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 73
$ws.Range("E21").Value = -31.506849315068
$ws.Range("F21").Value = 276
$ws.Range("G21").Value = 312
$ws.Range("H21").Value = -11.538461538461
$ws.Range("I21").Value = 3299
$ws.Range("J21").Value = 3495
$ws.Range("K21").Value = -5.608011444921
$ws.Range("L21").Value = -12.562947256824
$ws.Range("M21").Value = 16.613644397313
$ws.Range("N21").Value = -64.644732611724
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 47
$ws.Range("K22").Value = -36.170212765957
$ws.Range("L22").Value = -60
$ws.Range("M22").Value = -52.380952380952
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -17.857142857142
$ws.Range("I23").Value = 257
$ws.Range("J23").Value = 305
$ws.Range("K23").Value = -15.737704918032
$ws.Range("L23").Value = -18.927444794952
$ws.Range("M23").Value = 21.800947867298
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 9.090909090909
$ws.Range("F24").Value = 185
$ws.Range("G24").Value = 149
$ws.Range("H24").Value = 24.161073825503
$ws.Range("I24").Value = 2186
$ws.Range("J24").Value = 2059
$ws.Range("K24").Value = 6.168042739193
$ws.Range("L24").Value = -16.945288753799
$ws.Range("M24").Value = 62.407132243685
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -46.153846153846
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = 8.620689655172
$ws.Range("I25").Value = 850
$ws.Range("J25").Value = 808
$ws.Range("K25").Value = 5.19801980198
$ws.Range("L25").Value = -29.693961952026
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = -31.03448275862
$ws.Range("F26").Value = 117
$ws.Range("G26").Value = 109
$ws.Range("H26").Value = 7.339449541284
$ws.Range("I26").Value = 1414
$ws.Range("J26").Value = 1149
$ws.Range("K26").Value = 23.063533507397
$ws.Range("L26").Value = 22.001725625539
$ws.Range("M26").Value = -15.481171548117
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -60
$ws.Range("I27").Value = 91
$ws.Range("J27").Value = 83
$ws.Range("K27").Value = 9.638554216867
$ws.Range("L27").Value = 28.169014084507
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 133
$ws.Range("J28").Value = 107
$ws.Range("K28").Value = 24.29906542056
$ws.Range("L28").Value = 23.148148148148
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -66.666666666666
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 9
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 52
$ws.Range("J29").Value = 59
$ws.Range("K29").Value = -11.864406779661
$ws.Range("L29").Value = -39.53488372093
$ws.Range("N29").Value = -81.944444444444
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -16.666666666666
$ws.Range("I30").Value = 50
$ws.Range("J30").Value = 50
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = -41.860465116279
$ws.Range("N30").Value = -80.842911877394
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = -20
$ws.Range("L33").Value = 0
